# Add column-width configuration to the sheet, plus new date / date-time
# columns (F: date, G: date2, H: dateTime), and refresh the sample data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Column widths: A-F and H = 8 characters, G (the new "date2" column)
#    is wider at 16 characters so the timestamp fits.
#    (COM ColumnWidth excludes the ~0.8333 char cell-margin that Excel
#    adds on save, so subtract it to land on the intended stored width.)
# ---------------------------------------------------------------------
$narrow = 7.166666666666667   # -> stored width 8.0
$wide   = 15.166666666666666  # -> stored width 16.0

for ($c = 1; $c -le 8; $c++) {
    if ($c -eq 7) {
        $ws.Columns.Item($c).ColumnWidth = $wide
    } else {
        $ws.Columns.Item($c).ColumnWidth = $narrow
    }
}

# ---------------------------------------------------------------------
# 2. New header cells for the extra columns, styled like the existing
#    header row (copy formats only, so we inherit the exact header
#    font/fill/border instead of rebuilding it by hand).
# ---------------------------------------------------------------------
$ws.Range("F1").Value = "date"
$ws.Range("G1").Value = "date2"
$ws.Range("H1").Value = "dateTime"
$ws.Range("A1").Copy()
$ws.Range("F1:H1").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# 3. New date / date-time columns, same number format + fill as the
#    other numeric data columns.
# ---------------------------------------------------------------------
$ws.Range("F2:F13").NumberFormat = "yyyy-mm-dd"
$ws.Range("G2:G13").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("H2:H13").NumberFormat = "yyyy-mm-dd hh:mm:ss"
$ws.Range("F2:H13").Interior.ColorIndex = 40

$dateVal = 45020.0
$dateTimeVal = 45020.68039681713

# ---------------------------------------------------------------------
# 4. Refresh the sample data (columns C, D, E) and populate the new
#    F, G, H columns for every data row.
# ---------------------------------------------------------------------
$rows = @{
    2  = @(946554.0,    915719.8125,     198620.0)
    3  = @(159781.0,    865246.5,        778848.0)
    4  = @(375230.0,    927798.1875,     438320.0)
    5  = @(809286.0,    313665.15625,    247158.0)
    6  = @(502536.0,    427313.6875,     297902.0)
    7  = @(755012.0,    9438.6337890625, 359439.0)
    8  = @(277095.0,    198209.9375,     742489.0)
    9  = @(833938.0,    904344.6875,     784051.0)
    10 = @(428110.0,    763257.5,        552668.0)
    11 = @(183983.0,    813212.0625,     24319.0)
    12 = @(652746.0,    247624.34375,    164234.0)
    13 = @(326721.0,    254885.375,      944602.0)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 3).Value = $vals[0]   # C - cost
    $ws.Cells.Item($r, 4).Value = $vals[1]   # D - cost2
    $ws.Cells.Item($r, 5).Value = $vals[2]   # E - long value
    $ws.Cells.Item($r, 6).Value = $dateVal       # F - date
    $ws.Cells.Item($r, 7).Value = $dateTimeVal   # G - date2
    $ws.Cells.Item($r, 8).Value = $dateTimeVal   # H - dateTime
}

# ---------------------------------------------------------------------
# 5. Keep the original frozen header row (row 1) intact.
# ---------------------------------------------------------------------
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
